$wb = $excel.ActiveWorkbook

# Helper: force a range's contents to be stored as literal TEXT (not
# auto-coerced to a number) without leaving any NumberFormat / style
# residue behind. A plain `Range.Value = "123"` assignment infers a
# number from a numeric-looking string (same as typing it into Excel),
# which isn't what we want here since the source data keeps these
# columns as text. Writing a text formula ("=""123""") and immediately
# collapsing it to its value with Copy/PasteSpecial(xlPasteValues)
# keeps the literal string without touching formatting.
function Set-TextValues($range, [string[]]$values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $cell = $range.Cells.Item(1, $i + 1)
        $escaped = $values[$i].Replace('"', '""')
        $cell.Formula = '="' + $escaped + '"'
    }
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# ------------------------------------------------------------------
# 1. Insert a new "2022-Q1" sheet before "总计", copying "2021-Q4" so
#    it inherits the same sheet-level properties / cell formatting,
#    then trim it down to the new data and rename it.
# ------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$zj = $wb.Worksheets.Item("总计")
$q4.Copy($zj)

# re-fetch "总计" — the handle captured before the copy doesn't track the
# collection shifting, so $zj.Index would still report the pre-copy slot.
$zj = $wb.Worksheets.Item("总计")
$new = $wb.Worksheets.Item($zj.Index - 1)
$new.Name = "2022-Q1"

# the copied sheet has 4 data rows (2021-Q4 had 4 funds); 2022-Q1 only
# has 2, so drop the extra two rows.
$new.Rows("4:5").Delete()

# header row (already has the right style/text from the copy, but set
# it explicitly so the content matches exactly)
Set-TextValues $new.Range("B1:H1") @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")

# data rows
$new.Range("A2").Value = 0
Set-TextValues $new.Range("B2:G2") @("159855","银华中证影视主题ETF","0.96","97.27","3.63","0.0348")
$new.Range("H2").Value = 10

$new.Range("A3").Value = 1
Set-TextValues $new.Range("B3:G3") @("516620","国泰中证影视主题ETF","0.33","96.08","3.56","0.0117")
$new.Range("H3").Value = 10

# ------------------------------------------------------------------
# 2. Update the "总计" sheet with a new top row for 2022-Q1 and shift
#    the existing quarters down.
# ------------------------------------------------------------------
$zj.Range("A2").Value = 0
Set-TextValues $zj.Range("B2") @("2022-Q1")
$zj.Range("C2").Value = 2
$zj.Range("D2").Value = 0.05

$zj.Range("A3").Value = 1
Set-TextValues $zj.Range("B3") @("2021-Q4")
$zj.Range("C3").Value = 4
$zj.Range("D3").Value = 0.22

$zj.Range("A4").Value = 2
Set-TextValues $zj.Range("B4") @("2021-Q3")
$zj.Range("C4").Value = 2
$zj.Range("D4").Value = 1.89

$zj.Range("A5").Value = 3
Set-TextValues $zj.Range("B5") @("2020-Q4")
$zj.Range("C5").Value = 1
$zj.Range("D5").Value = 0.09
